$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.194969979067764
$ws.Range("C2").Value = 0.3616021647164303
$ws.Range("E2").Value = 0.4290078762398224
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.1039218777233302
$ws.Range("H2").Value = 0.2774058782005184
$ws.Range("O2").Value = 0.6527845062938411
$ws.Range("B3").Value = 1.04350217468874
$ws.Range("C3").Value = 0.3229059361183602
$ws.Range("E3").Value = 0.374192887699067
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.1064006724716862
$ws.Range("H3").Value = 0.2833977127817562
$ws.Range("O3").Value = 0.670301836390081
$ws.Range("B4").Value = 0.9501355902031605
$ws.Range("C4").Value = 0.2990604683654965
$ws.Range("E4").Value = 0.3406221926653075
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.1082086931143209
$ws.Range("H4").Value = 0.2873613110186071
$ws.Range("O4").Value = 0.6822505827476988
$ws.Range("B5").Value = 0.9119986520393013
$ws.Range("C5").Value = 0.2893224091187108
$ws.Range("E5").Value = 0.326961735161106
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.1090166775702919
$ws.Range("H5").Value = 0.2890478265128245
$ws.Range("O5").Value = 0.6874177402280992
$ws.Range("B6").Value = 0.9056607210990251
$ws.Range("C6").Value = 0.2877041753627338
$ws.Range("E6").Value = 0.3246945764550446
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.1091551219579436
$ws.Range("H6").Value = 0.28933217132845
$ws.Range("O6").Value = 0.68829367788166
$ws.Range("B7").Value = 0.9496216204440771
$ws.Range("C7").Value = 0.2989292208691552
$ws.Range("E7").Value = 0.3404378847639151
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.1082193025428744
$ws.Range("H7").Value = 0.2873837675229858
$ws.Range("O7").Value = 0.6823190651927291
$ws.Range("B8").Value = 1.142820978135035
$ws.Range("C8").Value = 0.3482779502538733
$ws.Range("E8").Value = 0.4100887711352783
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.1047167658497585
$ws.Range("H8").Value = 0.2794126499412926
$ws.Range("O8").Value = 0.6585756352619256
$ws.Range("B9").Value = 1.518703176012252
$ws.Range("C9").Value = 0.4443419307663703
$ws.Range("E9").Value = 0.5474405348454354
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.1001502592159227
$ws.Range("H9").Value = 0.2660502090460568
$ws.Range("O9").Value = 0.6215727217789748
$ws.Range("B10").Value = 1.792957655765861
$ws.Range("C10").Value = 0.5144581232898986
$ws.Range("E10").Value = 0.6489470192608593
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.09824318585554437
$ws.Range("H10").Value = 0.2576306654744798
$ws.Range("O10").Value = 0.6003390621484073
$ws.Range("B11").Value = 1.917292658480505
$ws.Range("C11").Value = 0.5462497296653623
$ws.Range("E11").Value = 0.695282065323795
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.09769958442565496
$ws.Range("H11").Value = 0.2541070735535556
$ws.Range("O11").Value = 0.5919987842164858
$ws.Range("B12").Value = 1.964312074147301
$ws.Range("C12").Value = 0.5582727241524594
$ws.Range("E12").Value = 0.7128529437601827
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.09754109883147066
$ws.Range("H12").Value = 0.2528171223847764
$ws.Range("O12").Value = 0.589032431333834
$ws.Range("B13").Value = 1.954188467178142
$ws.Range("C13").Value = 0.5556840716496936
$ws.Range("E13").Value = 0.7090676070841937
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.09757311252323575
$ws.Range("H13").Value = 0.2530929591759303
$ws.Range("O13").Value = 0.5896627164504906
$ws.Range("B14").Value = 1.921162261815709
$ws.Range("C14").Value = 0.54723918953664
$ws.Range("E14").Value = 0.6967271240135204
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.09768559198045779
$ws.Range("H14").Value = 0.2540000579030846
$ws.Range("O14").Value = 0.5917508811016603
$ws.Range("B15").Value = 1.900924408324613
$ws.Range("C15").Value = 0.5420643759546806
$ws.Range("E15").Value = 0.6891715086487835
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.09776068038752328
$ws.Range("H15").Value = 0.2545614669927403
$ws.Range("O15").Value = 0.5930550040939977
$ws.Range("B16").Value = 1.784823311747118
$ws.Range("C16").Value = 0.5123782973763582
$ws.Range("E16").Value = 0.6459222815661576
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.09828529682335585
$ws.Range("H16").Value = 0.2578671321789656
$ws.Range("O16").Value = 0.6009108537839012
$ws.Range("B17").Value = 1.713488491734779
$ws.Range("C17").Value = 0.4941395178017842
$ws.Range("E17").Value = 0.6194324140812455
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.0986906567322805
$ws.Range("H17").Value = 0.2599737534740569
$ws.Range("O17").Value = 0.6060696004334716
$ws.Range("B18").Value = 1.672418757450259
$ws.Range("C18").Value = 0.4836392530671105
$ws.Range("E18").Value = 0.6042110160391445
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.09895425259016122
$ws.Range("H18").Value = 0.2612142587046549
$ws.Range("O18").Value = 0.6091607925311422
$ws.Range("B19").Value = 1.658506485837904
$ws.Range("C19").Value = 0.4800823836917516
$ws.Range("E19").Value = 0.599059808982048
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.09904870647839914
$ws.Range("H19").Value = 0.2616392158645908
$ws.Range("O19").Value = 0.6102286484810406
$ws.Range("B20").Value = 1.721086345670244
$ws.Range("C20").Value = 0.4960820869129066
$ws.Range("E20").Value = 0.622250752681424
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.09864434901289343
$ws.Range("H20").Value = 0.2597465139530897
$ws.Range("O20").Value = 0.6055075913620698
$ws.Range("B21").Value = 1.930864603206032
$ws.Range("C21").Value = 0.5497200908074547
$ws.Range("E21").Value = 0.7003511346284341
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.09765126233887855
$ws.Range("H21").Value = 0.2537324149402309
$ws.Range("O21").Value = 0.5911323096306091
$ws.Range("B22").Value = 2.067595272575943
$ws.Range("C22").Value = 0.5846832081238063
$ws.Range("E22").Value = 0.7515397135821189
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.09727866267655116
$ws.Range("H22").Value = 0.2500605293194624
$ws.Range("O22").Value = 0.5828569914403516
$ws.Range("B23").Value = 1.994654345970616
$ws.Range("C23").Value = 0.5660314374250106
$ws.Range("E23").Value = 0.7242054608978776
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.09745197311760023
$ws.Range("H23").Value = 0.2519965215287741
$ws.Range("O23").Value = 0.5871704806097
$ws.Range("B24").Value = 1.717651537621464
$ws.Range("C24").Value = 0.4952038966787882
$ws.Range("E24").Value = 0.6209765571265962
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.09866518963342941
$ws.Range("H24").Value = 0.2598491574396391
$ws.Range("O24").Value = 0.6057612853372802
$ws.Range("B25").Value = 1.417345556484349
$ws.Range("C25").Value = 0.418433092748387
$ws.Range("E25").Value = 0.5101881373765309
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.1011345695931851
$ws.Range("H25").Value = 0.2694206889789612
$ws.Range("O25").Value = 0.6305467712937229
